$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 520609
$ws.Range("I2").Value = 727652.7
$ws.Range("J2").Value = 2999.75
$ws.Range("K2").Value = 727652.7
$ws.Range("L2").Value = 2999.75
$ws.Range("M2").Value = -727539.7
$ws.Range("N2").Value = -3225.75
$ws.Range("H9").Value = 72.916664
$ws.Range("I9").Value = 70.5
$ws.Range("K9").Value = 70.5
$ws.Range("M9").Value = 98.5
$ws.Range("H33").Value = 3099.3242
$ws.Range("I33").Value = 3881.3928
$ws.Range("K33").Value = 3881.3928
$ws.Range("M33").Value = -3652.3928
$ws.Range("H64").Value = 5332.6665
$ws.Range("I64").Value = 5433.5
$ws.Range("K64").Value = 5433.5
$ws.Range("M64").Value = -5185.5
$ws.Range("H67").Value = 5332.6665
$ws.Range("I67").Value = 5433.5
$ws.Range("K67").Value = 5433.5
$ws.Range("M67").Value = -4575.5
$ws.Range("H69").Value = 15845
$ws.Range("J69").Value = 15845
$ws.Range("L69").Value = 47535
$ws.Range("N69").Value = -49283
$ws.Range("H72").Value = 15845
$ws.Range("J72").Value = 15845
$ws.Range("L72").Value = 142605
$ws.Range("N72").Value = -151341
$ws.Range("H113").Value = 6858.0835
$ws.Range("J113").Value = 7042
$ws.Range("L113").Value = 7042
$ws.Range("N113").Value = -13550
$ws.Range("H135").Value = 129231.875
$ws.Range("J135").Value = 205184.6
$ws.Range("L135").Value = 1846661.4
$ws.Range("N135").Value = -1851731.4

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12588.487
$ws.Range("I32").Value = 9029.75
$ws.Range("K32").Value = 9029.75
$ws.Range("M32").Value = -8742.75
$ws.Range("H45").Value = 252550.75
$ws.Range("I45").Value = 435629.9
$ws.Range("K45").Value = 435629.9
$ws.Range("M45").Value = -435252.9

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1628
$ws.Range("I20").Value = 1546.5555
$ws.Range("J20").Value = 1872.3334
$ws.Range("K20").Value = 1546.5555
$ws.Range("L20").Value = 1872.3334
$ws.Range("M20").Value = -1299.5555
$ws.Range("N20").Value = -2366.3334
$ws.Range("H107").Value = 1553.2727
$ws.Range("I107").Value = 1103.4
$ws.Range("K107").Value = 1103.4
$ws.Range("M107").Value = 816.5999999999999

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 9179.25
$ws.Range("I41").Value = 9179.25
$ws.Range("K41").Value = 9179.25
$ws.Range("M41").Value = -8751.25
$ws.Range("H107").Value = 111472.81
$ws.Range("I107").Value = 150367.34
$ws.Range("K107").Value = 150367.34
$ws.Range("M107").Value = -148447.34

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 14285746
$ws.Range("J2").Value = 40000010
$ws.Range("L2").Value = 240000060
$ws.Range("N2").Value = -240000286
$ws.Range("H11").Value = 468.7
$ws.Range("I11").Value = 468.7
$ws.Range("K11").Value = 1406.1
$ws.Range("M11").Value = -1266.1
$ws.Range("H69").Value = 963.38464
$ws.Range("J69").Value = 928.25
$ws.Range("L69").Value = 2784.75
$ws.Range("N69").Value = -4406.75
$ws.Range("H72").Value = 963.38464
$ws.Range("J72").Value = 928.25
$ws.Range("L72").Value = 8354.25
$ws.Range("N72").Value = -16466.25
$ws.Range("H98").Value = 1482
$ws.Range("I98").Value = 479
$ws.Range("J98").Value = 1768.5714
$ws.Range("K98").Value = 1437
$ws.Range("L98").Value = 5305.7142
$ws.Range("M98").Value = 61
$ws.Range("N98").Value = -8301.7142
$ws.Range("H113").Value = 1871.2
$ws.Range("I113").Value = 1369.1666
$ws.Range("J113").Value = 2624.25
$ws.Range("K113").Value = 4107.4998
$ws.Range("L113").Value = 7872.75
$ws.Range("M113").Value = -1937.4998
$ws.Range("N113").Value = -12212.75

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5899.5
$ws.Range("H73").Value = 5899.5
$ws.Range("H80").Value = 3714
$ws.Range("J80").Value = 4921.5
$ws.Range("L80").Value = 4921.5
$ws.Range("N80").Value = -6917.5
$ws.Range("H83").Value = 3714
$ws.Range("J83").Value = 4921.5
$ws.Range("L83").Value = 24607.5
$ws.Range("N83").Value = -34591.5
$ws.Range("H132").Value = 34613.78
$ws.Range("I132").Value = 44432.375
$ws.Range("J132").Value = 5158
$ws.Range("K132").Value = 133297.125
$ws.Range("L132").Value = 15474
$ws.Range("M132").Value = -130767.125
$ws.Range("N132").Value = -20534
$ws.Range("H139").Value = 65997.5
$ws.Range("J139").Value = 65997.5
$ws.Range("L139").Value = 65997.5
$ws.Range("N139").Value = -76277.5

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 917552.0600000001
$ws.Range("I7").Value = 917552.0600000001
$ws.Range("K7").Value = 917552.0600000001
$ws.Range("M7").Value = -917440.0600000001
$ws.Range("H14").Value = 9999.5
$ws.Range("I14").Value = 9999
$ws.Range("K14").Value = 9999
$ws.Range("M14").Value = -9827
$ws.Range("H46").Value = 11522.615
$ws.Range("I46").Value = 46500
$ws.Range("J46").Value = 5163.091
$ws.Range("K46").Value = 46500
$ws.Range("L46").Value = 5163.091
$ws.Range("M46").Value = -46312
$ws.Range("N46").Value = -5539.091
$ws.Range("H61").Value = 1907.2693
$ws.Range("I61").Value = 1907.2693
$ws.Range("K61").Value = 1907.2693
$ws.Range("M61").Value = -1705.2693
$ws.Range("H113").Value = 1907.2693
$ws.Range("I113").Value = 1907.2693
$ws.Range("K113").Value = 1907.2693
$ws.Range("M113").Value = 262.7307000000001
$ws.Range("H126").Value = 917552.0600000001
$ws.Range("I126").Value = 917552.0600000001
$ws.Range("K126").Value = 2752656.18
$ws.Range("M126").Value = -2750186.18
$ws.Range("H132").Value = 90805.86
$ws.Range("I132").Value = 154248.38
$ws.Range("K132").Value = 462745.14
$ws.Range("M132").Value = -460215.14

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 46160.332
$ws.Range("I132").Value = 45777.16
$ws.Range("J132").Value = 50950
$ws.Range("K132").Value = 137331.48
$ws.Range("L132").Value = 152850
$ws.Range("M132").Value = -134801.48
$ws.Range("N132").Value = -157910
$ws.Range("H140").Value = 110000
$ws.Range("J140").Value = 110000
$ws.Range("L140").Value = 110000
$ws.Range("N140").Value = -120360
